$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column B
$ws.Range("B2").Value = 88
$ws.Range("B3").Value = -259
$ws.Range("B5").Value = 0.343
$ws.Range("B6").Value = -0.618
$ws.Range("B7").Value = 0.343
$ws.Range("B8").Value = -0.618

# Move the active selection from B7 to B9
$ws.Range("B9").Select()
